# Auto-generated edit script: refresh market-price derived columns (H-N)
# across the 8 job sheets, per the scheduled-runner data update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 6263.625
$ws.Range("I18").Value = 2935.0833
$ws.Range("K18").Value = 2935.0833
$ws.Range("M18").Value = -2651.0833

# Row 19
$ws.Range("H19").Value = 1768
$ws.Range("I19").Value = 1713.5714
$ws.Range("K19").Value = 1713.5714
$ws.Range("M19").Value = -1538.5714

# Row 100
$ws.Range("H100").Value = 1931.9231
$ws.Range("I100").Value = 2265
$ws.Range("J100").Value = 100
$ws.Range("K100").Value = 2265
$ws.Range("L100").Value = 100
$ws.Range("M100").Value = -1724
$ws.Range("N100").Value = -1182

# Row 112
$ws.Range("H112").Value = 6085.9355
$ws.Range("J112").Value = 6340.1724
$ws.Range("L112").Value = 19020.5172
$ws.Range("N112").Value = -21236.5172

# Row 116
$ws.Range("H116").Value = 971471.1
$ws.Range("I116").Value = 1271524.8
$ws.Range("K116").Value = 1271524.8
$ws.Range("M116").Value = -1268082.8

# Row 127
$ws.Range("H127").Value = 1220.1428
$ws.Range("I127").Value = 741.6
$ws.Range("K127").Value = 2224.8
$ws.Range("M127").Value = 2735.2

# Row 129
$ws.Range("H129").Value = 1022.93335
$ws.Range("I129").Value = 667.53845
$ws.Range("J129").Value = 3333
$ws.Range("K129").Value = 2002.61535
$ws.Range("L129").Value = 9999
$ws.Range("M129").Value = 2997.38465
$ws.Range("N129").Value = -19999

# Row 132
$ws.Range("H132").Value = 34932.81
$ws.Range("I132").Value = 38557.684
$ws.Range("K132").Value = 115673.052
$ws.Range("M132").Value = -113143.052

# Row 137
$ws.Range("H137").Value = 20433.658
$ws.Range("I137").Value = 52050.168
$ws.Range("J137").Value = 14505.5625
$ws.Range("K137").Value = 156150.504
$ws.Range("L137").Value = 43516.6875
$ws.Range("M137").Value = -153600.504
$ws.Range("N137").Value = -48616.6875

# Row 138
$ws.Range("H138").Value = 43141.08
$ws.Range("J138").Value = 114688.445
$ws.Range("L138").Value = 344065.335
$ws.Range("N138").Value = -354345.335

# Row 141
$ws.Range("H141").Value = 1503.1875
$ws.Range("I141").Value = 1470.0667
$ws.Range("K141").Value = 4410.2001
$ws.Range("M141").Value = 769.7999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22069.059
$ws.Range("I32").Value = 24062.361
$ws.Range("J32").Value = 3332
$ws.Range("K32").Value = 24062.361
$ws.Range("L32").Value = 3332
$ws.Range("M32").Value = -23775.361
$ws.Range("N32").Value = -3906

# Row 45
$ws.Range("H45").Value = 4672.263
$ws.Range("I45").Value = 3998.5
$ws.Range("J45").Value = 4851.933
$ws.Range("K45").Value = 3998.5
$ws.Range("L45").Value = 4851.933
$ws.Range("M45").Value = -3621.5
$ws.Range("N45").Value = -5605.933

# Row 61
$ws.Range("H61").Value = 6006.0347
$ws.Range("I61").Value = 806.7273
$ws.Range("J61").Value = 22346.715
$ws.Range("K61").Value = 806.7273
$ws.Range("L61").Value = 22346.715
$ws.Range("M61").Value = -594.7273
$ws.Range("N61").Value = -22770.715

# Row 63
$ws.Range("H63").Value = 3159.6
$ws.Range("I63").Value = 1519.2
$ws.Range("J63").Value = 4800
$ws.Range("K63").Value = 1519.2
$ws.Range("L63").Value = 4800
$ws.Range("M63").Value = -833.2
$ws.Range("N63").Value = -6172

# Row 66
$ws.Range("H66").Value = 3159.6
$ws.Range("I66").Value = 1519.2
$ws.Range("J66").Value = 4800
$ws.Range("K66").Value = 7596
$ws.Range("L66").Value = 24000
$ws.Range("M66").Value = -4164
$ws.Range("N66").Value = -30864

# Row 80
$ws.Range("H80").Value = 29989.666

# Row 83
$ws.Range("H83").Value = 29989.666

# Row 92
$ws.Range("H92").Value = 60000
$ws.Range("J92").Value = 60000
$ws.Range("L92").Value = 60000
$ws.Range("N92").Value = -64992

# Row 122
$ws.Range("H122").Value = 3065.889
$ws.Range("I122").Value = 2941.8572
$ws.Range("K122").Value = 8825.571599999999
$ws.Range("M122").Value = -6375.571599999999

# Row 134
$ws.Range("H134").Value = 80000
$ws.Range("I134").Value = 80000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 80000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -74930
$ws.Range("N134").ClearContents()

# Row 136
$ws.Range("H136").Value = 6006.0347
$ws.Range("I136").Value = 806.7273
$ws.Range("J136").Value = 22346.715
$ws.Range("K136").Value = 2420.1819
$ws.Range("L136").Value = 67040.145
$ws.Range("M136").Value = 129.8181
$ws.Range("N136").Value = -72140.145

# Row 139
$ws.Range("H139").Value = 116674.375
$ws.Range("J139").Value = 116674.375
$ws.Range("L139").Value = 116674.375
$ws.Range("N139").Value = -126954.375

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 14020

# Row 85
$ws.Range("H85").Value = 14020

# Row 105
$ws.Range("H105").Value = 2056.4
$ws.Range("I105").Value = 1861.9048
$ws.Range("J105").Value = 3077.5
$ws.Range("K105").Value = 1861.9048
$ws.Range("L105").Value = 3077.5
$ws.Range("M105").Value = -114.9048
$ws.Range("N105").Value = -6571.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5560961
$ws.Range("I31").Value = 7147950
$ws.Range("K31").Value = 7147950
$ws.Range("M31").Value = -7147655

# Row 34
$ws.Range("H34").Value = 5560961
$ws.Range("I34").Value = 7147950
$ws.Range("K34").Value = 7147950
$ws.Range("M34").Value = -7147748

# Row 94
$ws.Range("H94").Value = 1602.3914
$ws.Range("I94").Value = 1129.4166
$ws.Range("J94").Value = 2118.3635
$ws.Range("K94").Value = 1129.4166
$ws.Range("L94").Value = 2118.3635
$ws.Range("M94").Value = -678.4166
$ws.Range("N94").Value = -3020.3635

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# Row 105
$ws.Range("H105").Value = 2474.6667
$ws.Range("J105").Value = 3999.25
$ws.Range("L105").Value = 3999.25
$ws.Range("N105").Value = -7493.25

# Row 122
$ws.Range("H122").Value = 1925.8948
$ws.Range("I122").Value = 1842.2142
$ws.Range("K122").Value = 5526.642599999999
$ws.Range("M122").Value = -3076.642599999999

# Row 131
$ws.Range("H131").Value = 65000
$ws.Range("J131").Value = 65000
$ws.Range("L131").Value = 65000
$ws.Range("N131").Value = -75080

$ws = $wb.Worksheets.Item("CUL")
# Row 114
$ws.Range("H114").Value = 1115.25
$ws.Range("J114").Value = 1800
$ws.Range("L114").Value = 5400
$ws.Range("N114").Value = -11908

# Row 122
$ws.Range("H122").Value = 960.2222
$ws.Range("J122").Value = 1020.2857
$ws.Range("L122").Value = 9182.5713
$ws.Range("N122").Value = -14082.5713

# Row 123
$ws.Range("H123").Value = 5157.5
$ws.Range("I123").Value = 4943.3335
$ws.Range("K123").Value = 14830.0005
$ws.Range("M123").Value = -12380.0005

# Row 129
$ws.Range("H129").Value = 2046.579
$ws.Range("I129").Value = 1782.3636
$ws.Range("K129").Value = 5347.0908
$ws.Range("M129").Value = -347.0907999999999

# Row 131
$ws.Range("H131").Value = 4344
$ws.Range("J131").Value = 2999
$ws.Range("L131").Value = 8997
$ws.Range("N131").Value = -19077

# Row 139
$ws.Range("H139").Value = 5248.467
$ws.Range("I139").Value = 5248.467
$ws.Range("K139").Value = 15745.401
$ws.Range("M139").Value = -10605.401

$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 50000
$ws.Range("J51").Value = 50000
$ws.Range("L51").Value = 50000
$ws.Range("N51").Value = -51018

# Row 132
$ws.Range("H132").Value = 3170.0715
$ws.Range("I132").Value = 2312.25
$ws.Range("J132").Value = 3513.2
$ws.Range("K132").Value = 6936.75
$ws.Range("L132").Value = 10539.6
$ws.Range("M132").Value = -4406.75
$ws.Range("N132").Value = -15599.6

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2580.611
$ws.Range("I93").Value = 2448.3
$ws.Range("J93").Value = 2746
$ws.Range("K93").Value = 2448.3
$ws.Range("L93").Value = 2746
$ws.Range("M93").Value = -1200.3
$ws.Range("N93").Value = -5242

# Row 100
$ws.Range("H100").Value = 2848.6
$ws.Range("I100").Value = 2754.5715
$ws.Range("K100").Value = 2754.5715
$ws.Range("M100").Value = -2213.5715

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 14047.167
$ws.Range("J45").Value = 14047.167
$ws.Range("L45").Value = 14047.167
$ws.Range("N45").Value = -15029.167

# Row 113
$ws.Range("H113").Value = 1014.65515
$ws.Range("I113").Value = 991.6316
$ws.Range("J113").Value = 1058.4
$ws.Range("K113").Value = 2974.8948
$ws.Range("L113").Value = 3175.2
$ws.Range("M113").Value = -804.8948
$ws.Range("N113").Value = -7515.200000000001

# Row 132
$ws.Range("H132").Value = 73158.8
$ws.Range("I132").Value = 178000
$ws.Range("J132").Value = 3264.6667
$ws.Range("K132").Value = 534000
$ws.Range("L132").Value = 9794.000100000001
$ws.Range("M132").Value = -531470
$ws.Range("N132").Value = -14854.0001
